# Update generated output numbers ("想去人数" column F) on the
# "展览" and "全部类型" sheets, rows 3-7.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 2133
    $ws.Range("F4").Value = 1619
    $ws.Range("F5").Value = 320
    $ws.Range("F6").Value = 1039
    $ws.Range("F7").Value = 509
}
